# Generate Report for Handback
# Update the "generate date" / "handoff/handback datetime" timestamps to reflect
# a newly regenerated handback report.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the first file row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-18 13:07:06"

# zh-cn sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime" for the first file row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-18 13:06:56"
$wsZhCn.Range("K2").Value = "2016-08-18 13:07:28"

# de-de sheet: "Correspond Handback DateTime" for the first file row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-18 13:07:35"
